# Auto-generated Excel COM-interop script
# Refreshes the Leve Profit market-data snapshot values (columns H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match the latest
# pull from the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 230.88235
$ws.Range("I12").Value = 226.5625
$ws.Range("K12").Value = 226.5625
$ws.Range("M12").Value = -56.5625
$ws.Range("H43").Value = 23822
$ws.Range("J43").Value = 51850.75
$ws.Range("L43").Value = 51850.75
$ws.Range("N43").Value = -51988.75
$ws.Range("H76").Value = 2801.6667
$ws.Range("I76").Value = 2801.6667
$ws.Range("K76").Value = 2801.6667
$ws.Range("M76").Value = -2486.6667
$ws.Range("H79").Value = 2801.6667
$ws.Range("I79").Value = 2801.6667
$ws.Range("K79").Value = 2801.6667
$ws.Range("M79").Value = -1709.6667
$ws.Range("H98").Value = 572.8182
$ws.Range("I98").Value = 572.8182
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 572.8182
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 925.1818
$ws.Range("N98").ClearContents()
$ws.Range("H116").Value = 6224.375
$ws.Range("J116").Value = 6174
$ws.Range("L116").Value = 6174
$ws.Range("N116").Value = -13058
$ws.Range("H122").Value = 572.8182
$ws.Range("I122").Value = 572.8182
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1718.4546
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 731.5454
$ws.Range("N122").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1638.5834
$ws.Range("I97").Value = 766.3
$ws.Range("K97").Value = 766.3
$ws.Range("M97").Value = -270.3
$ws.Range("H102").Value = 6035.4546
$ws.Range("I102").Value = 5298.75
$ws.Range("K102").Value = 5298.75
$ws.Range("M102").Value = -3676.75
$ws.Range("H122").Value = 3571.9092
$ws.Range("I122").Value = 3413.1428
$ws.Range("K122").Value = 10239.4284
$ws.Range("M122").Value = -7789.428400000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2459
$ws.Range("I20").Value = 2459
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 2459
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -2212
$ws.Range("N20").ClearContents()
$ws.Range("H80").Value = 1351.5
$ws.Range("J80").Value = 1475.1666
$ws.Range("L80").Value = 1475.1666
$ws.Range("N80").Value = -3471.1666
$ws.Range("H83").Value = 1351.5
$ws.Range("J83").Value = 1475.1666
$ws.Range("L83").Value = 7375.833000000001
$ws.Range("N83").Value = -17359.833
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H105").Value = 2963.5833
$ws.Range("I105").Value = 2819.0908
$ws.Range("K105").Value = 2819.0908
$ws.Range("M105").Value = -1072.0908
$ws.Range("H134").Value = 2427.95
$ws.Range("I134").Value = 2632.7693
$ws.Range("K134").Value = 7898.3079
$ws.Range("M134").Value = -5363.3079

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6828
$ws.Range("I86").Value = 6234.6665
$ws.Range("K86").Value = 6234.6665
$ws.Range("M86").Value = -5111.6665
$ws.Range("H89").Value = 6828
$ws.Range("I89").Value = 6234.6665
$ws.Range("K89").Value = 31173.3325
$ws.Range("M89").Value = -25557.3325
$ws.Range("H99").Value = 15418.931
$ws.Range("I99").Value = 11808.789
$ws.Range("J99").Value = 22278.2
$ws.Range("K99").Value = 11808.789
$ws.Range("L99").Value = 22278.2
$ws.Range("M99").Value = -10310.789
$ws.Range("N99").Value = -25274.2
$ws.Range("H126").Value = 15418.931
$ws.Range("I126").Value = 11808.789
$ws.Range("J126").Value = 22278.2
$ws.Range("K126").Value = 35426.367
$ws.Range("L126").Value = 66834.60000000001
$ws.Range("M126").Value = -32956.367
$ws.Range("N126").Value = -71774.60000000001
$ws.Range("H132").Value = 2925.7144
$ws.Range("J132").Value = 5776.5
$ws.Range("L132").Value = 17329.5
$ws.Range("N132").Value = -22389.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 2500
$ws.Range("I48").Value = 2500
$ws.Range("K48").Value = 7500
$ws.Range("M48").Value = -7250
$ws.Range("H68").Value = 1015.4286
$ws.Range("I68").Value = 1027.75
$ws.Range("J68").Value = 999
$ws.Range("K68").Value = 3083.25
$ws.Range("L68").Value = 2997
$ws.Range("M68").Value = -2272.25
$ws.Range("N68").Value = -4619
$ws.Range("H71").Value = 1015.4286
$ws.Range("I71").Value = 1027.75
$ws.Range("J71").Value = 999
$ws.Range("K71").Value = 9249.75
$ws.Range("L71").Value = 8991
$ws.Range("M71").Value = -5193.75
$ws.Range("N71").Value = -17103

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 12630
$ws.Range("I55").Value = 12630
$ws.Range("K55").Value = 12630
$ws.Range("M55").Value = -12303
$ws.Range("H97").Value = 889
$ws.Range("I97").Value = 884
$ws.Range("K97").Value = 884
$ws.Range("M97").Value = -388
$ws.Range("H132").Value = 3317.8333
$ws.Range("I132").Value = 2681.5
$ws.Range("K132").Value = 8044.5
$ws.Range("M132").Value = -5514.5
$ws.Range("H138").Value = 66500
$ws.Range("J138").Value = 66500
$ws.Range("L138").Value = 66500
$ws.Range("N138").Value = -76780

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4194
$ws.Range("I62").Value = 4194
$ws.Range("K62").Value = 4194
$ws.Range("M62").Value = -3570
$ws.Range("H65").Value = 4194
$ws.Range("I65").Value = 4194
$ws.Range("K65").Value = 20970
$ws.Range("M65").Value = -17850
$ws.Range("H96").Value = 760419
$ws.Range("I96").Value = 1104564.9
$ws.Range("J96").Value = 3298
$ws.Range("K96").Value = 1104564.9
$ws.Range("L96").Value = 3298
$ws.Range("M96").Value = -1103191.9
$ws.Range("N96").Value = -6044
$ws.Range("H100").Value = 1506.4
$ws.Range("I100").Value = 1324.4286
$ws.Range("K100").Value = 2648.8572
$ws.Range("M100").Value = -2107.8572
$ws.Range("H132").Value = 7959.04
$ws.Range("I132").Value = 5764.2
$ws.Range("J132").Value = 16738.4
$ws.Range("K132").Value = 17292.6
$ws.Range("L132").Value = 50215.2
$ws.Range("M132").Value = -14762.6
$ws.Range("N132").Value = -55275.2

Write-Host "Applied 158 value updates and 5 cell clears across 8 sheets."
